# Clean up before editing the paper
#
# Column A (a study-level "hyper/hypo/hyper-hypo" grouping tag) was only
# populated on the first data row of each contiguous group. This fills the
# tag down onto every other data row of the same group so every row in the
# sheet carries its own value (separator rows that have no "B" entry are
# intentionally left alone). It also corrects row 494, which previously
# pointed at the stray "hyper(cont)" label, to use the standard "hyper"
# label instead - once nothing references "hyper(cont)" any more it drops
# out of the shared-strings table automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3:A15").Value = "hyper"
$ws.Range("A17:A22").Value = "hyper"
$ws.Range("A25:A29").Value = "hyper"
$ws.Range("A31:A36").Value = "hyper"
$ws.Range("A39:A44").Value = "hyper"
$ws.Range("A47").Value = "hyper"
$ws.Range("A50:A57").Value = "hyper"
$ws.Range("A62:A66").Value = "hyper"
$ws.Range("A69:A81").Value = "hyper"
$ws.Range("A85:A95").Value = "hyper"
$ws.Range("A99:A109").Value = "hyper"
$ws.Range("A112:A128").Value = "hyper"
$ws.Range("A131:A149").Value = "hyper"
$ws.Range("A152:A180").Value = "hyper"
$ws.Range("A183:A186").Value = "hypo"
$ws.Range("A189:A191").Value = "hyper"
$ws.Range("A194:A202").Value = "hyper"
$ws.Range("A205:A209").Value = "hyper/hypo"
$ws.Range("A212:A218").Value = "hyper/hypo"
$ws.Range("A224:A227").Value = "hyper"
$ws.Range("A229:A243").Value = "hyper/hypo"
$ws.Range("A247").Value = "hyper"
$ws.Range("A250:A256").Value = "hyper/hypo"
$ws.Range("A259:A260").Value = "hypo"
$ws.Range("A262:A271").Value = "hypo"
$ws.Range("A274:A278").Value = "hypo"
$ws.Range("A281:A285").Value = "hypo"
$ws.Range("A289").Value = "hypo"
$ws.Range("A292").Value = "hypo"
$ws.Range("A294:A304").Value = "hypo"
$ws.Range("A307:A315").Value = "hypo"
$ws.Range("A319:A366").Value = "hypo"
$ws.Range("A369:A407").Value = "hypo"
$ws.Range("A410:A434").Value = "hypo"
$ws.Range("A437:A451").Value = "hypo"
$ws.Range("A454:A482").Value = "hypo"
$ws.Range("A485:A492").Value = "hypo"
$ws.Range("A494:A507").Value = "hyper"

# Leave the cursor where it ends up after this cleanup pass (just past the
# last data row), matching the saved view state.
$ws.Range("A518").Select()
